# Update the raw per-trial timing measurements (rows 2-7, columns B-K) on the
# "Brute Force" and "Divide and Conquer" sheets with re-measured values,
# including a fixed base case (n=1, row 2) which now reports 0 on both
# sheets. The Mean / Std Dev columns (L/M) already hold AVERAGE / STDEV.S
# formulas over B:K, so they recalculate automatically once the raw values
# are replaced - no formula changes are needed.

$wb = $excel.ActiveWorkbook

# Rows correspond to input sizes 1, 200, 400, 600, 800, 1000 (sheet rows 2-7).
# Columns correspond to Trial 1..10 (sheet columns B-K).
$bruteForceData = @(
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0.05807828903198242, 0.05686426162719727, 0.05745244026184082, 0.0593106746673584, 0.05231833457946777, 0.05089449882507324, 0.07591700553894043, 0.04929590225219727, 0.05038046836853027, 0.05597138404846191),
    @(0.2040259838104248, 0.2096037864685059, 0.2194952964782715, 0.2248122692108154, 0.2093229293823242, 0.221970796585083, 0.2260551452636719, 0.2072842121124268, 0.2137517929077148, 0.2056353092193604),
    @(0.4634144306182861, 0.4614531993865967, 0.6120820045471191, 0.4511904716491699, 0.4643123149871826, 0.475156307220459, 0.5719192028045654, 0.4635899066925049, 0.4749021530151367, 0.4656960964202881),
    @(0.8650112152099609, 0.82700514793396, 0.9095418453216553, 0.9233608245849609, 0.8982982635498047, 0.9550714492797852, 0.8983595371246338, 0.9080491065979004, 0.9102048873901367, 0.9253726005554199),
    @(1.434290409088135, 1.415933132171631, 1.590587615966797, 1.633319854736328, 1.539464712142944, 1.50659441947937, 1.43195104598999, 1.413079023361206, 1.379208087921143, 1.407622098922729)
)

$divideAndConquerData = @(
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0.005215883255004883, 0.008206605911254883, 0, 0.00806427001953125, 0.006873846054077148, 0, 0),
    @(0, 0, 0.01570534706115723, 0.009578466415405273, 0, 0, 0.003987550735473633, 0, 0, 0),
    @(0.0157933235168457, 0, 0, 0, 0, 0.01604771614074707, 0.004022598266601562, 0.0101161003112793, 0, 0),
    @(0.01618456840515137, 0, 0.01617574691772461, 0.003621101379394531, 0.009146213531494141, 0, 0, 0.01602816581726074, 0.009388208389282227, 0.01594018936157227),
    @(0.01586008071899414, 0.01311087608337402, 0.01558566093444824, 0.01580691337585449, 0.01010394096374512, 0.01583266258239746, 0.01607418060302734, 0.01607418060302734, 0, 0.008038997650146484)
)

function Set-TrialData($worksheet, $data) {
    for ($rowIdx = 0; $rowIdx -lt $data.Length; $rowIdx++) {
        $sheetRow = $rowIdx + 2
        $rowValues = $data[$rowIdx]
        for ($colIdx = 0; $colIdx -lt $rowValues.Length; $colIdx++) {
            $sheetCol = $colIdx + 2
            $worksheet.Cells.Item($sheetRow, $sheetCol).Value = $rowValues[$colIdx]
        }
    }
}

$wsBrute = $wb.Worksheets.Item("Brute Force")
Set-TrialData $wsBrute $bruteForceData

$wsDivide = $wb.Worksheets.Item("Divide and Conquer")
Set-TrialData $wsDivide $divideAndConquerData
